# Apply the new bug-report row (row 9) and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "look command has no message for empty location"
$ws.Range("C9").Value = "ricky"
$ws.Range("D9").Value = "empty space location has no long description"
$ws.Range("E9").Value = "added various long descriptions with random selection"
$ws.Range("F9").Value = "fixed"

$ws.Range("D7").Select()
